# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the newer scrape snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1546
$ws1.Range("F8").Value  = 171
$ws1.Range("F9").Value  = 751
$ws1.Range("F12").Value = 337
$ws1.Range("F13").Value = 59
$ws1.Range("F16").Value = 23
$ws1.Range("F21").Value = 15468
$ws1.Range("F22").Value = 1536
$ws1.Range("F26").Value = 11096
$ws1.Range("F27").Value = 761

# Sheet "全部类型": row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1546
$ws4.Range("F9").Value  = 171
$ws4.Range("F10").Value = 751
$ws4.Range("F14").Value = 337
$ws4.Range("F15").Value = 59
$ws4.Range("F19").Value = 23
$ws4.Range("F24").Value = 15468
$ws4.Range("F25").Value = 1536
$ws4.Range("F29").Value = 11096
$ws4.Range("F30").Value = 761
